$p = $ppt.ActivePresentation

# "Next Steps" slide (Lesson 8.1): the first bullet of the content
# placeholder referenced two example files. Remove the reference to
# 08-1-decode.rkt (and the "and 08-2-" connector), leaving only a
# reference to 08-1-merge-sort.rkt.
$oldSpan = "files 08-1-decode.rkt and 08-2-merge-sort"
$newSpan = "files 08-1-merge-sort"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi)
                $idx = $para.Text.IndexOf($oldSpan)
                if ($idx -ge 0) {
                    $sel = $para.Characters($idx + 1, $oldSpan.Length)
                    $sel.Text = $newSpan
                }
            }
        }
    }
}
